# atualiza script do grid
# - rename the "Planilha1" sheet (sheetId 9) to "tesouro"
# - make that sheet ("tesouro") the active/selected tab instead of "grid_export"
# - set its selection to AL14 (was the full A1:BA42 block)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Planilha1")
$ws.Name = "tesouro"

# Activating this sheet updates workbookView.activeTab and moves
# tabSelected from "grid_export" onto "tesouro" automatically.
$ws.Activate()

# Move the selection/active cell to AL14 on the now-active sheet.
$ws.Range("AL14").Select()
